$wb = $excel.ActiveWorkbook

# Reference cell far outside the used range, with pristine default styling,
# used to strip any auto-applied column styling/number-format from newly
# written cells so they keep the default (unstyled) cell format.

$active = $wb.Worksheets.Item("Active")
$activeRef = $active.Range("Z1000")
$active.Range("A2").Value = "Castle-Play-Tent-Orian"
$active.Range("B2").Value = "ORIAN Princess Castle Playhouse Tent for Girls with LED Star Lights – Indoor & Outdoor Large Kids Play Tent for Imaginative Games – ASTM Certified, Princess Tent, 230 Polyester Taffeta. Pink 55""x53""."
$active.Range("C2").Value = "B07TV8HCMT"
$active.Range("D2").Value = "'793611000049"
$active.Range("E2").Value = "missing"
$active.Range("F2").Value = "'10.78"
$active.Range("A2:F2").Style = $activeRef.Style
$active.Range("A3").Value = "GA-DBZG-GTH3"
$active.Range("B3").Value = "Orian Toys Teepee Tent for Kids: Child’s Indoor Outdoor Canvas Fairytale Tipi Playroom, LED Star Lights, Easy Assembly, 59 by 45 Inches, Ages 3+"
$active.Range("C3").Value = "B08N56JXNR"
$active.Range("D3").Value = "B08N56JXNR"
$active.Range("E3").Value = "child"
$active.Range("F3").Value = "'13.6"
$active.Range("A3:F3").Style = $activeRef.Style
$active.Range("A4").Value = "STEM-Toys-Education-Engineering-Orian"
$active.Range("B4").Value = "Orian Toys 5 in 1 STEM Learning Toys for Boys and Girls, Best IQ Builder STEM Learning Toys Creative Construction Engineering for Kids 5-11 years old, DIY Building Kit, 132 Pieces, Play Set - Gift Box"
$active.Range("C4").Value = "B08B1P25HR"
$active.Range("D4").Value = "'793611000094"
$active.Range("E4").Value = "missing"
$active.Range("F4").Value = "'5.05"
$active.Range("A4:F4").Style = $activeRef.Style
$active.Range("A5").Value = "Teepee-Tent-Orian"
$active.Range("B5").Value = "Teepee Tent for Kids - A Fairytale Tipi Tent Kids Love. LED Star Lights, Dream Catcher - Strong Indoor Tee Pee Tent - Kids Play Tent for Boys & Girls"
$active.Range("C5").Value = "B089M859ZQ"
$active.Range("D5").Value = "'793611000087"
$active.Range("E5").Value = "child"
$active.Range("F5").Value = "'13.6"
$active.Range("A5:F5").Style = $activeRef.Style
$active.Range("A6").Value = "Teepee_Pompon"
$active.Range("B6").Value = "Orian Teepee Tent for Kids - Indoor Tent, Pompon Ball Design, Built-in Mat, Inner Pockets, Window, LED Star Lights, Dream Catcher, Unique Threaded Poles - The Perfect Play Tent for Girls & Boys"
$active.Range("C6").Value = "B00KG50JUU"
$active.Range("D6").Value = "B00KG50JUU"
$active.Range("E6").Value = "child"
$active.Range("F6").Value = "'14.1"
$active.Range("A6:F6").Style = $activeRef.Style
$active.Range("A7").Value = "XQ-OECU-MGA4"
$active.Range("B7").Value = "Orian Pop Beads Jewelry Making Kit for Girls, 550+ Piece Set, Pop Beads for Girls Ages 3 and Up, Fun and Colorful Snap Beads, Bracelet Making Kit, Necklaces and Rings. Great Colorful Unicorn Gift Bag"
$active.Range("C7").Value = "B07T8RLTY2"
$active.Range("D7").Value = "'793611000032"
$active.Range("E7").Value = "missing"
$active.Range("F7").Value = "'5.7"
$active.Range("A7:F7").Style = $activeRef.Style

$incomplete = $wb.Worksheets.Item("InComplete")
$incompleteRef = $incomplete.Range("Z1000")
$incomplete.Range("A2").Value = "Teepee Tents"
$incomplete.Range("B2").Value = "Orian Teepee Tent for Kids - Indoor Tent Built-in Mat, Inner Pockets, Window, LED Star Lights, Dream Catcher - The Perfect Play Tent for Girls & Boys"
$incomplete.Range("C2").Value = "B0B136RGLH"
$incomplete.Range("D2").Value = "B0B136RGLH"
$incomplete.Range("E2").Value = "parent"
$incomplete.Range("F2").Value = "'13.6"
$incomplete.Range("A2:F2").Style = $incompleteRef.Style

